$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 538) holds a "Förändrad" (changed/updated) date
# that was bumped by one day (2023-09-08 -> 2023-09-09, serial 45177 -> 45178).
$ws.Range("C2:C538").Value = 45178
